$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A6").Value = "Ajudante"
